$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bus")

$k1 = $ws.Range("K1")
$k1.Value = "X"
$k1.Font.Bold = $true
$k1.HorizontalAlignment = -4108
$k1.VerticalAlignment = -4160
$k1.Borders.Item(7).LineStyle = 1
$k1.Borders.Item(10).LineStyle = 1

$l1 = $ws.Range("L1")
$l1.Value = "Y"
$l1.Font.Bold = $true
$l1.HorizontalAlignment = -4108
$l1.VerticalAlignment = -4160
$l1.Borders.Item(7).LineStyle = 1
$l1.Borders.Item(10).LineStyle = 1

$m1 = $ws.Range("M1")
$m1.Value = "storage"
$m1.Font.Bold = $true
$m1.HorizontalAlignment = -4108
$m1.VerticalAlignment = -4160
$m1.Borders.Item(7).LineStyle = 1
$m1.Borders.Item(10).LineStyle = 1

$dummy = $ws.Range("Z100")
$dummy.Value = 1
$dummy.Value = 1
